$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 135
$ws.Range("H135").Value = 1353.0186
$ws.Range("I135").Value = 1164.775
$ws.Range("J135").Value = 1890.8572
$ws.Range("K135").Value = 10482.975
$ws.Range("L135").Value = 17017.7148
$ws.Range("M135").Value = -7947.975
$ws.Range("N135").Value = -22087.7148

# Row 137
$ws.Range("H137").Value = 4167485.8
$ws.Range("I137").Value = 784.6129
$ws.Range("J137").Value = 11765588
$ws.Range("K137").Value = 2353.8387
$ws.Range("L137").Value = 35296764
$ws.Range("M137").Value = 196.1613000000002
$ws.Range("N137").Value = -35301864

# Row 141
$ws.Range("H141").Value = 1236.475
$ws.Range("I141").Value = 955.1177
$ws.Range("J141").Value = 2830.8333
$ws.Range("K141").Value = 2865.3531
$ws.Range("L141").Value = 8492.499899999999
$ws.Range("M141").Value = 2314.6469
$ws.Range("N141").Value = -18852.4999

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 3702
$ws.Range("I45").Value = 10012
$ws.Range("J45").Value = 2440
$ws.Range("K45").Value = 10012
$ws.Range("L45").Value = 2440
$ws.Range("M45").Value = -9635
$ws.Range("N45").Value = -3194

# Row 61
$ws.Range("H61").Value = 7813655
$ws.Range("I61").Value = 10205172
$ws.Range("J61").Value = 1364.3334
$ws.Range("K61").Value = 10205172
$ws.Range("L61").Value = 1364.3334
$ws.Range("M61").Value = -10204960
$ws.Range("N61").Value = -1788.3334

# Row 136
$ws.Range("H136").Value = 7813655
$ws.Range("I136").Value = 10205172
$ws.Range("J136").Value = 1364.3334
$ws.Range("K136").Value = 30615516
$ws.Range("L136").Value = 4093.0002
$ws.Range("M136").Value = -30612966
$ws.Range("N136").Value = -9193.0002

$ws = $wb.Worksheets.Item("BSM")
# Row 107
$ws.Range("H107").Value = 1072.2963
$ws.Range("I107").Value = 997.7143
$ws.Range("J107").Value = 1333.3334
$ws.Range("K107").Value = 997.7143
$ws.Range("L107").Value = 1333.3334
$ws.Range("M107").Value = 922.2857
$ws.Range("N107").Value = -5173.3334

# Row 108
$ws.Range("H108").Value = 40684
$ws.Range("J108").Value = 40684
$ws.Range("L108").Value = 40684
$ws.Range("N108").Value = -48364

# Row 134
$ws.Range("H134").Value = 1855.2885
$ws.Range("I134").Value = 1161.1915
$ws.Range("J134").Value = 8379.799999999999
$ws.Range("K134").Value = 3483.5745
$ws.Range("L134").Value = 25139.4
$ws.Range("M134").Value = -948.5744999999997
$ws.Range("N134").Value = -30209.4

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 5379325.5
$ws.Range("I31").Value = 3354.0889
$ws.Range("J31").Value = 19609838
$ws.Range("K31").Value = 3354.0889
$ws.Range("L31").Value = 19609838
$ws.Range("M31").Value = -3059.0889
$ws.Range("N31").Value = -19610428

# Row 34
$ws.Range("H34").Value = 5379325.5
$ws.Range("I34").Value = 3354.0889
$ws.Range("J34").Value = 19609838
$ws.Range("K34").Value = 3354.0889
$ws.Range("L34").Value = 19609838
$ws.Range("M34").Value = -3152.0889
$ws.Range("N34").Value = -19610242

# Row 134
$ws.Range("H134").Value = 1382.5526
$ws.Range("I134").Value = 1525.72
$ws.Range("J134").Value = 1107.2307
$ws.Range("K134").Value = 4577.16
$ws.Range("L134").Value = 3321.6921
$ws.Range("M134").Value = -2042.16
$ws.Range("N134").Value = -8391.6921

# Row 140
$ws.Range("H140").Value = 44794.145
$ws.Range("J140").Value = 44794.145
$ws.Range("L140").Value = 44794.145
$ws.Range("N140").Value = -55154.145

$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 5000166
$ws.Range("I4").Value = 5000166
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 15000498
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -15000386
$ws.Range("N4").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Range("H122").Value = 2780409.8
$ws.Range("I122").Value = 3510939.2
$ws.Range("J122").Value = 4397.4
$ws.Range("K122").Value = 10532817.6
$ws.Range("L122").Value = 13192.2
$ws.Range("M122").Value = -10530367.6
$ws.Range("N122").Value = -18092.2

# Row 132
$ws.Range("H132").Value = 3291.157
$ws.Range("I132").Value = 2349.3076
$ws.Range("J132").Value = 6352.1665
$ws.Range("K132").Value = 7047.9228
$ws.Range("L132").Value = 19056.4995
$ws.Range("M132").Value = -4517.9228
$ws.Range("N132").Value = -24116.4995

# Row 134
$ws.Range("H134").Value = 38930.3
$ws.Range("J134").Value = 38930.3
$ws.Range("L134").Value = 116790.9
$ws.Range("N134").Value = -121860.9

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 4713.97
$ws.Range("I7").Value = 4474.9565
$ws.Range("J7").Value = 5237.524
$ws.Range("K7").Value = 4474.9565
$ws.Range("L7").Value = 5237.524
$ws.Range("M7").Value = -4362.9565
$ws.Range("N7").Value = -5461.524

# Row 22
$ws.Range("H22").Value = 999.1070999999999
$ws.Range("I22").Value = 820
$ws.Range("J22").Value = 1028.9584
$ws.Range("K22").Value = 820
$ws.Range("L22").Value = 1028.9584
$ws.Range("M22").Value = -525
$ws.Range("N22").Value = -1618.9584

# Row 27
$ws.Range("H27").Value = 999.1070999999999
$ws.Range("I27").Value = 820
$ws.Range("J27").Value = 1028.9584
$ws.Range("K27").Value = 820
$ws.Range("L27").Value = 1028.9584
$ws.Range("M27").Value = -713
$ws.Range("N27").Value = -1242.9584

# Row 122
$ws.Range("H122").Value = 4913.6665
$ws.Range("I122").Value = 4634.8887
$ws.Range("J122").Value = 5750
$ws.Range("K122").Value = 13904.6661
$ws.Range("L122").Value = 17250
$ws.Range("M122").Value = -11454.6661
$ws.Range("N122").Value = -22150

# Row 126
$ws.Range("H126").Value = 4713.97
$ws.Range("I126").Value = 4474.9565
$ws.Range("J126").Value = 5237.524
$ws.Range("K126").Value = 13424.8695
$ws.Range("L126").Value = 15712.572
$ws.Range("M126").Value = -10954.8695
$ws.Range("N126").Value = -20652.572

# Row 139
$ws.Range("H139").Value = 60663.57
$ws.Range("J139").Value = 60663.57
$ws.Range("L139").Value = 60663.57
$ws.Range("N139").Value = -70943.57000000001

$ws = $wb.Worksheets.Item("WVR")
# Row 2
$ws.Range("H2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()

# Row 132
$ws.Range("H132").Value = 1301.5763
$ws.Range("I132").Value = 929.9375
$ws.Range("J132").Value = 2923.2727
$ws.Range("K132").Value = 2789.8125
$ws.Range("L132").Value = 8769.8181
$ws.Range("M132").Value = -259.8125
$ws.Range("N132").Value = -13829.8181
